$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pthlh"
$ws.Range("C2").Value = "Pth1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1203763333333333
$ws.Range("H2").Value = 0.361129
$ws.Range("I2").Value = 0.007481788818141147
$ws.Range("J2").Value = 0.007481788818141146
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.278930666666667
$ws.Range("N2").Value = 3.836792
$ws.Range("O2").Value = 0.2016523671602801
$ws.Range("P2").Value = 0.2016523671602801
$ws.Range("Q2").Value = 0.1539529842408889
$ws.Range("R2").Value = 1.385576858168
$ws.Range("S2").Value = 0.001508720425771477
$ws.Range("T2").Value = 0.001508720425771477

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pthlh"
$ws.Range("C3").Value = "Pth1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1203763333333333
$ws.Range("H3").Value = 0.361129
$ws.Range("I3").Value = 0.007481788818141147
$ws.Range("J3").Value = 0.007481788818141146
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.645054333333333
$ws.Range("N3").Value = 4.935163
$ws.Range("O3").Value = 0.2593800501230792
$ws.Range("P3").Value = 0.2593800501230792
$ws.Range("Q3").Value = 0.1980256087807778
$ws.Range("R3").Value = 1.782230479027
$ws.Range("S3").Value = 0.001940626758659745
$ws.Range("T3").Value = 0.001940626758659744

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pthlh"
$ws.Range("C4").Value = "Pth1r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1203763333333333
$ws.Range("H4").Value = 0.361129
$ws.Range("I4").Value = 0.007481788818141147
$ws.Range("J4").Value = 0.007481788818141146
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.418269666666667
$ws.Range("N4").Value = 10.254809
$ws.Range("O4").Value = 0.5389675827166406
$ws.Range("P4").Value = 0.5389675827166406
$ws.Range("Q4").Value = 0.411478768817889
$ws.Range("R4").Value = 3.703308919361
$ws.Range("S4").Value = 0.004032441633709925
$ws.Range("T4").Value = 0.004032441633709924

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pthlh"
$ws.Range("C5").Value = "Pth1r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.12404233333333
$ws.Range("H5").Value = 42.372127
$ws.Range("I5").Value = 0.8778561289441074
$ws.Range("J5").Value = 0.8778561289441073
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.278930666666667
$ws.Range("N5").Value = 3.836792
$ws.Range("O5").Value = 0.2016523671602801
$ws.Range("P5").Value = 0.2016523671602801
$ws.Range("Q5").Value = 18.06367087739822
$ws.Range("R5").Value = 162.573037896584
$ws.Range("S5").Value = 0.1770217664277393
$ws.Range("T5").Value = 0.1770217664277393

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pthlh"
$ws.Range("C6").Value = "Pth1r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.12404233333333
$ws.Range("H6").Value = 42.372127
$ws.Range("I6").Value = 0.8778561289441074
$ws.Range("J6").Value = 0.8778561289441073
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.645054333333333
$ws.Range("N6").Value = 4.935163
$ws.Range("O6").Value = 0.2593800501230792
$ws.Range("P6").Value = 0.2593800501230792
$ws.Range("Q6").Value = 23.23481704463344
$ws.Range("R6").Value = 209.113353401701
$ws.Range("S6").Value = 0.2276983667263749
$ws.Range("T6").Value = 0.2276983667263749

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pthlh"
$ws.Range("C7").Value = "Pth1r"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.12404233333333
$ws.Range("H7").Value = 42.372127
$ws.Range("I7").Value = 0.8778561289441074
$ws.Range("J7").Value = 0.8778561289441073
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.418269666666667
$ws.Range("N7").Value = 10.254809
$ws.Range("O7").Value = 0.5389675827166406
$ws.Range("P7").Value = 0.5389675827166406
$ws.Range("Q7").Value = 48.27978547874923
$ws.Range("R7").Value = 434.518069308743
$ws.Range("S7").Value = 0.4731359957899931
$ws.Range("T7").Value = 0.473135995789993

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Pthlh"
$ws.Range("C8").Value = "Pth1r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.844826333333334
$ws.Range("H8").Value = 5.534479000000001
$ws.Range("I8").Value = 0.1146620822377516
$ws.Range("J8").Value = 0.1146620822377516
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.278930666666667
$ws.Range("N8").Value = 3.836792
$ws.Range("O8").Value = 0.2016523671602801
$ws.Range("P8").Value = 0.2016523671602801
$ws.Range("Q8").Value = 2.359404972374223
$ws.Range("R8").Value = 21.23464475136801
$ws.Range("S8").Value = 0.02312188030676932
$ws.Range("T8").Value = 0.02312188030676932

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Pthlh"
$ws.Range("C9").Value = "Pth1r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.844826333333334
$ws.Range("H9").Value = 5.534479000000001
$ws.Range("I9").Value = 0.1146620822377516
$ws.Range("J9").Value = 0.1146620822377516
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.645054333333333
$ws.Range("N9").Value = 4.935163
$ws.Range("O9").Value = 0.2593800501230792
$ws.Range("P9").Value = 0.2593800501230792
$ws.Range("Q9").Value = 3.034839553897445
$ws.Range("R9").Value = 27.31355598507701
$ws.Range("S9").Value = 0.02974105663804465
$ws.Range("T9").Value = 0.02974105663804464

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Pthlh"
$ws.Range("C10").Value = "Pth1r"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.844826333333334
$ws.Range("H10").Value = 5.534479000000001
$ws.Range("I10").Value = 0.1146620822377516
$ws.Range("J10").Value = 0.1146620822377516
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.418269666666667
$ws.Range("N10").Value = 10.254809
$ws.Range("O10").Value = 0.5389675827166406
$ws.Range("P10").Value = 0.5389675827166406
$ws.Range("Q10").Value = 6.306113895501223
$ws.Range("R10").Value = 56.75502505951101
$ws.Range("S10").Value = 0.06179914529293763
$ws.Range("T10").Value = 0.06179914529293763
